$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.945.26"
$ws.Range("E2").Value = "  +1.50%  "

$ws.Range("D3").Value = "3.202.79"
$ws.Range("E3").Value = "  +0.98%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'602.90"
$ws.Range("E5").Value = "  +4.01%  "

$ws.Range("D6").Value = "'152.50"
$ws.Range("E6").Value = "  +0.53%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "3.201.01"
$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("E10").Value = "  -1.48%  "

$ws.Range("D11").Value = "'6.16"
$ws.Range("E11").Value = "  -1.56%  "

$ws.Range("D12").Value = "'0.505"
$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("D13").Value = "'0.0000270"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("D14").Value = "'38.26"
$ws.Range("E14").Value = "  +1.50%  "

$ws.Range("D15").Value = "3.727.50"
$ws.Range("E15").Value = "  +0.93%  "

$ws.Range("D16").Value = "66.057.05"
$ws.Range("E16").Value = "  +1.51%  "

$ws.Range("D17").Value = "'7.42"
$ws.Range("E17").Value = "  +3.61%  "

$ws.Range("D18").Value = "3.213.89"
$ws.Range("E18").Value = "  +1.23%  "

$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("D20").Value = "'508.43"
$ws.Range("E20").Value = "  -1.07%  "

$ws.Range("D21").Value = "'15.46"
$ws.Range("E21").Value = "  +4.10%  "

$ws.Range("D22").Value = "'0.729"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Value = "'15.20"
$ws.Range("E23").Value = "  -0.45%  "

$ws.Range("D24").Value = "'7.98"
$ws.Range("E24").Value = "  +2.41%  "

$ws.Range("D25").Value = "'84.98"
$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  +2.58%  "

$ws.Range("D28").Value = "'9.10"
$ws.Range("E28").Value = "  +1.36%  "

$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  +2.34%  "

$ws.Range("D30").Value = "'2.84"
$ws.Range("E30").Value = "  +2.12%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'27.99"
$ws.Range("E31").Value = "  +0.38%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.75"
$ws.Range("E32").Value = "  +6.77%  "

$ws.Range("E33").Value = "  +0.15%  "

$ws.Range("E34").Value = "  +1.05%  "

$ws.Range("D35").Value = "'6.57"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("E36").Value = "  -0.85%  "

$ws.Range("D37").Value = "'0.0900"
$ws.Range("E37").Value = "  +0.44%  "

$ws.Range("D38").Value = "'477.62"
$ws.Range("E38").Value = "  +0.82%  "

$ws.Range("D39").Value = "'0.0418"
$ws.Range("E39").Value = "  -0.64%  "

$ws.Range("D40").Value = "'2.94"
$ws.Range("E40").Value = "  -6.30%  "

$ws.Range("D41").Value = "'8.81"
$ws.Range("E41").Value = "  +1.78%  "

$ws.Range("D42").Value = "'0.294"
$ws.Range("E42").Value = "  +2.79%  "

$ws.Range("D43").Value = "'0.118"
$ws.Range("E43").Value = "  -0.77%  "

$ws.Range("D44").Value = "2.929.30"
$ws.Range("E44").Value = "  -4.49%  "

$ws.Range("D45").Value = "'2.43"
$ws.Range("E45").Value = "  +1.01%  "

$ws.Range("E46").Value = "  +4.42%  "

$ws.Range("D47").Value = "'28.55"
$ws.Range("E47").Value = "  -1.73%  "

$ws.Range("D49").Value = "'0.115"
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("D50").Value = "'2.29"
$ws.Range("E50").Value = "  +1.68%  "

$ws.Range("D51").Value = "'33.85"
$ws.Range("E51").Value = "  +4.05%  "
